{"js": "// The document has several \"<id>...</id>\" labels that were each split\n// across three separate runs: the opening \"<id>\" tag, the inner id\n// value, and the closing \"</id>\" tag. This collapses each one back into\n// a single run (keeping the \"<id>\" run's formatting), and for the first\n// occurrence also corrects the id value itself (p030r_a2 -> p030r_2);\n// the other four keep their existing id value, just merged into one run.\n\nconst body = context.document.body;\n\nconst replacements = [\n  { oldText: \"<id>p030r_a2</id>\", newText: \"<id>p030r_2</id>\" },\n  { oldText: \"<id>p031r_1</id>\", newText: \"<id>p031r_1</id>\" },\n  { oldText: \"<id>p031r_2</id>\", newText: \"<id>p031r_2</id>\" },\n  { oldText: \"<id>p031r_3</id>\", newText: \"<id>p031r_3</id>\" },\n  { oldText: \"<id>p031r_4</id>\", newText: \"<id>p031r_4</id>\" },\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains several \"<id>...</id>\" runs that were each split\n# across three separate w:r runs (open tag / inner value / close tag),\n# e.g. \"<id>\" + \"p030r_a2\" + \"</id>\". This consolidates each occurrence\n# into a single run (keeping the formatting of the opening \"<id>\" run)\n# and, for the first one, also updates the id value itself\n# (p030r_a2 -> p030r_2); the other four keep the same id value.\n\n$doc = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"<id>p030r_a2</id>\"; New = \"<id>p030r_2</id>\" },\n    @{ Old = \"<id>p031r_1</id>\";  New = \"<id>p031r_1</id>\" },\n    @{ Old = \"<id>p031r_2</id>\";  New = \"<id>p031r_2</id>\" },\n    @{ Old = \"<id>p031r_3</id>\";  New = \"<id>p031r_3</id>\" },\n    @{ Old = \"<id>p031r_4</id>\";  New = \"<id>p031r_4</id>\" }\n)\n\nforeach ($item in $replacements) {\n    $range = $doc.Content\n    $find = $range.Find\n    $find.Text = $item.Old\n    $find.Replacement.Text = $item.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
